$wb = $excel.ActiveWorkbook

# Replace "Ready for handoff" with "In Translation" wherever it appears
# across all worksheets (Overview, zh-cn, de-de). This is a localization
# status report being regenerated for archive: entries that were
# previously marked ready for handoff are now back "In Translation".
foreach ($ws in $wb.Worksheets) {
    $used = $ws.UsedRange
    foreach ($cell in $used.Cells) {
        # NOTE: keep the literal on the left of -eq; Cells can hold
        # booleans (e.g. "True"/"False" values) and PowerShell coerces
        # the comparison to the type of the left-hand operand, so a
        # boolean cell would otherwise match any non-empty string.
        if ("Ready for handoff" -eq $cell.Value2) {
            $cell.Value = "In Translation"
        }
    }
}

# Narrow the status columns to fit the shorter text:
#   Overview sheet  -> columns E (zh-cn) and F (de-de)
#   zh-cn / de-de sheets -> column C (Status)
$overview = $wb.Worksheets.Item("Overview")
$overview.Columns.Item(5).ColumnWidth = 12.5
$overview.Columns.Item(6).ColumnWidth = 12.5

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Columns.Item(3).ColumnWidth = 12.5

$dede = $wb.Worksheets.Item("de-de")
$dede.Columns.Item(3).ColumnWidth = 12.5
